$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "CS001;CS004"
$ws.Range("E9").Select()
